$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "nombre"
$ws.Range("B1").Value = "usuario"
$ws.Range("C1").Value = "correo"
$ws.Range("D1").Value = "password"
$ws.Range("E1").Value = "perfil"

# --- Row 2: Mario Monreal ---
$ws.Range("A2").Value = "Mario Monreal"
$ws.Range("B2").Value = "moma2"
$ws.Range("C2").Value = "pema@gloablhitss.com"
$ws.Range("D2").Value = "p3r4l35#4"
$ws.Range("E2").Value = 59

# --- Row 3: Martina Flores ---
$ws.Range("A3").Value = "Martina Flores"
$ws.Range("B3").Value = "marti1"
$ws.Range("C3").Value = "marti1@globalhits.com"
$ws.Range("D3").Value = "c4rr3ra#1"
$ws.Range("E3").Value = 59

# --- Row 4: Jorge Casio (new hyperlink required on C4) ---
$ws.Range("A4").Value = "Jorge Casio"
$ws.Range("B4").Value = "casjo86"
$ws.Range("C4").Value = "casjo21@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:casjo21@gmail.com")
$ws.Range("C4").Style = "Hipervínculo"
$ws.Range("D4").Value = "c45j0123"
$ws.Range("E4").Value = 59

# --- Selection moves to A5 ---
$ws.Range("A5").Select()
